$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.150.11"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.791.99"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "224.42"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "32.75"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").Value = "0.0707"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "0.0931"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "2.049.80"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "1.804.52"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "10.89"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").Value = "0.626"
$ws.Range("E15").Value = "  -2.60%  "
$ws.Range("D16").Value = "34.150.14"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "4.18"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "68.03"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "243.50"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "10.74"
$ws.Range("E22").Value = "  -4.01%  "
$ws.Range("E23").Value = "  -3.31%  "
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("D25").Value = "159.51"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("D26").Value = "16.31"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").Value = "7.07"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").Value = "3.52"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("D35").Value = "1.394.73"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").Value = "0.647"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").Value = "2.22"
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").Value = "78.99"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("D42").Value = "0.918"
$ws.Range("E42").Value = "  -4.87%  "
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").Value = "0.0₆0149"
$ws.Range("E44").Value = "  +19.42%  "
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "108.13"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "5.90"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "1.949.31"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").Value = "12.21"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  +0.04%  "
